# "Added last minute updates"
#
# The document's first paragraph holds the topic placeholder
# "**ID__AFFARS_5344_topic_7__ID**" followed by a single trailing-space
# run. This edit:
#   1. gives that paragraph a paragraph border (5pt space on all sides),
#   2. widens its left indent from 120 twips (6pt) to 225 twips (11.25pt),
#   3. drops the stray trailing-space run, and
#   4. renames the placeholder to "**ID__AFFARS_5344_302__ID**".

$d = $word.ActiveDocument

# The placeholder paragraph is the very first paragraph in the document.
$p1 = $d.Paragraphs(1)
$pf = $p1.Range.ParagraphFormat

# --- pPr: add <w:pBdr> (top/left/bottom/right, w:space="5") ------------
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5

# --- pPr: <w:ind w:left="120"/> -> <w:ind w:left="225"/> ---------------
# ParagraphFormat.LeftIndent is in points; 225 twips = 11.25 points.
$pf.LeftIndent = 11.25

# --- locate the placeholder run -----------------------------------------
# Find.Execute (without a replacement) leaves the range collapsed onto the
# matched text, so its End gives us the exact boundary right after the
# placeholder, independent of any hard-coded character offsets.
$findRange = $d.Content
$found = $findRange.Find.Execute("**ID__AFFARS_5344_topic_7__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$placeholderEnd = $findRange.End

# --- remove the trailing " " run that follows the placeholder ----------
$spaceRunRange = $d.Range($placeholderEnd, $placeholderEnd + 1)
$spaceRunRange.Delete()

# --- rename the placeholder text ----------------------------------------
$d.Content.Find.Execute("**ID__AFFARS_5344_topic_7__ID**", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_5344_302__ID**", 2)
